$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A24").Value = 23
$ws.Range("B24").Value = 26
$ws.Range("C24").Value = 0
$ws.Range("D24").Value = 25
$ws.Range("E24").Value = 39
$ws.Range("F24").Value = 51
$ws.Range("G24").Value = 90

$ws.Range("A25").Value = 24
$ws.Range("B25").Value = 26
$ws.Range("C25").Value = 1
$ws.Range("D25").Value = 27
$ws.Range("E25").Value = 36
$ws.Range("F25").Value = 54
$ws.Range("G25").Value = 90
